# Correcting some documentation issues:
# - the worksheet was still carrying the legacy "H21R00" board name; rename
#   it to match this module (H1FR50).
# - leave the cursor/selection where the author left it when the file was
#   last saved (D6:F6) instead of the old D4:F4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "H1FR50"

$ws.Range("D6:F6").Select()
